$d = $word.ActiveDocument

# The document has a set of runs whose rPr currently contains only a <w:b/>
# toggle. Every such run needs two additional explicit toggles appended
# right after <w:b/>: <w:i w:val="false"/> and <w:strike w:val="false"/>.
# The single run that also carries a <w:color/> (the "www.j2eestar.com"
# hyperlink-ish text inside the "Cell Phone..." paragraph) must be left
# untouched, so that paragraph is handled as a special case below using
# Find to scope the operation to only the preceding run's text.

$specialParaText = "Cell Phone: (916) 812-1709 Email: spektr44@hotmail.com URL: "

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    $pStart = $p.Range.Start
    $pEnd = $p.Range.End

    if ($pEnd -le $pStart) {
        continue
    }

    # Exclude the trailing paragraph-mark character from the range so we
    # never touch the paragraph mark's own run properties.
    $rangeEnd = $pEnd - 1
    if ($rangeEnd -le $pStart) {
        continue
    }

    $fullText = $d.Range($pStart, $rangeEnd).Text

    if ($fullText -eq $specialParaText + "www.j2eestar.com") {
        # Only the first run ("Cell Phone: ... URL: ") gets the new
        # toggles; the trailing "www.j2eestar.com" run (has <w:color/>)
        # is left alone.
        $scoped = $d.Range($pStart, $rangeEnd)
        $found = $scoped.Find.Execute($specialParaText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
        if ($found) {
            $scoped.Font.Italic = 0
            $scoped.Font.StrikeThrough = 0
        }
        continue
    }

    $r = $d.Range($pStart, $rangeEnd)
    $r.Font.Italic = 0
    $r.Font.StrikeThrough = 0
}
